$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

foreach ($r in 2..3) {
    $ws.Cells.Item($r, 4).Value = 0.0174      # D
    $ws.Cells.Item($r, 5).Value = 0.0809      # E
    $ws.Cells.Item($r, 6).Value = -0.062      # F

    $ws.Cells.Item($r, 9).Value = 0           # I
    $ws.Cells.Item($r, 10).Value = 0          # J
    $ws.Cells.Item($r, 11).Value = 159.7      # K
    $ws.Cells.Item($r, 12).Value = 0.281757233592096   # L
    $ws.Cells.Item($r, 13).Value = 0.182      # M
    $ws.Cells.Item($r, 14).Value = 0.0001624274877286925  # N
    $ws.Cells.Item($r, 15).Value = 0.001139636819035692   # O
    $ws.Cells.Item($r, 16).Value = 0.182      # P
    $ws.Cells.Item($r, 17).Value = 0.0001624274877286925  # Q
    $ws.Cells.Item($r, 18).Value = 0.001139636819035692   # R

    $ws.Cells.Item($r, 21).Value = 700.6      # U
    $ws.Cells.Item($r, 22).Value = 0.625256581883088      # V
    $ws.Cells.Item($r, 23).Value = 0.0881492520836783     # W
    $ws.Cells.Item($r, 24).Value = 0.05188135250274772    # X
    $ws.Cells.Item($r, 25).Value = 0.03626789958093057    # Y
    $ws.Cells.Item($r, 26).Value = 0.3639983302828886     # Z
    $ws.Cells.Item($r, 27).Value = 0          # AA
    $ws.Cells.Item($r, 28).Value = 0.04466536796585908    # AB
    $ws.Cells.Item($r, 29).Value = -0.04466536796585908   # AC
    $ws.Cells.Item($r, 30).Value = 706        # AD
    $ws.Cells.Item($r, 31).Value = 0          # AE
    $ws.Cells.Item($r, 32).Value = 706        # AF
    $ws.Cells.Item($r, 33).Value = 5.399999999999977      # AG
    $ws.Cells.Item($r, 34).Value = 0.3865316178483438     # AH
    $ws.Cells.Item($r, 35).Value = 0.2483467004361896     # AI
    $ws.Cells.Item($r, 36).Value = 0.004796163069544344   # AJ
    $ws.Cells.Item($r, 37).Value = 0.002520773037064689   # AK

    # AN (40) and AP (42) are removed (cleared)
    $ws.Cells.Item($r, 40).ClearContents()
    $ws.Cells.Item($r, 42).ClearContents()
}
